$wb = $excel.ActiveWorkbook

# ============================================================
# Overview sheet
# ============================================================
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = "f36151f6-bf51-4878-9d26-ad0381d1a250.md"
$ov.Range("G2").Value = "2016-08-30 17:10:50"

# Update the B2 hyperlink display text, keeping its original target URL.
$ov.Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4e801975cb7803105c37071f9fea236d961090bb/e2e/bdaeab86-0f22-48e3-bccb-b33bcdc08d2b.md", [Type]::Missing, [Type]::Missing, "e2e\f36151f6-bf51-4878-9d26-ad0381d1a250.md")

# ============================================================
# zh-cn sheet
# ============================================================
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("G2").Value = "f36151f6-bf51-4878-9d26-ad0381d1a250.b560f3ffa54a3f09b9545d5fbb167cc772ba37ac.zh-cn.xlf"
$zh.Range("H2").Value = "2016-08-30 17:10:45"
$zh.Range("J2").Value = ""
$zh.Range("K2").Value = "0001-01-01 00:00:00"

# Clear "Latest Target File" (I2) back to a blank, unstyled cell and drop its
# hyperlink; only the A2 (source file name) hyperlink/display survives.
$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4e801975cb7803105c37071f9fea236d961090bb/e2e/bdaeab86-0f22-48e3-bccb-b33bcdc08d2b.md", [Type]::Missing, [Type]::Missing, "f36151f6-bf51-4878-9d26-ad0381d1a250.md")
$zh.Range("I2").Value = ""
$zh.Range("I2").Style = "Normal"

$zh.Columns.Item(9).ColumnWidth = 18.6506053379604
$zh.Columns.Item(10).ColumnWidth = 21.7054770333426

# ============================================================
# de-de sheet
# ============================================================
$de = $wb.Worksheets.Item("de-de")

$de.Range("G2").Value = "f36151f6-bf51-4878-9d26-ad0381d1a250.b560f3ffa54a3f09b9545d5fbb167cc772ba37ac.de-de.xlf"
$de.Range("H2").Value = "2016-08-30 17:10:50"
$de.Range("J2").Value = ""
$de.Range("K2").Value = "0001-01-01 00:00:00"

$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4e801975cb7803105c37071f9fea236d961090bb/e2e/bdaeab86-0f22-48e3-bccb-b33bcdc08d2b.md", [Type]::Missing, [Type]::Missing, "f36151f6-bf51-4878-9d26-ad0381d1a250.md")
$de.Range("I2").Value = ""
$de.Range("I2").Style = "Normal"

$de.Columns.Item(9).ColumnWidth = 18.6506053379604
$de.Columns.Item(10).ColumnWidth = 21.7054770333426
